# Update crypto price/volume data on the active worksheet.
# Mirrors the refreshed values from the "Updated cryptos list" GitHub Actions run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.981.00"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +1.38%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.308.79"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +5.68%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "598.50"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.25%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.06"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +4.34%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.304.77"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +5.84%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.61%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +2.07%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +3.28%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +2.36%  "
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.66%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.82"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +1.15%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.851.62"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +5.74%  "
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +1.10%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.307.53"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +5.77%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "64.057.20"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +1.39%  "
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +1.95%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "481.17"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.91%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.29"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.49%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.744"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +6.11%  "
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +3.99%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.49"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +3.25%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "84.35"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -3.60%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +2.21%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +1.98%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.26"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +3.18%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +2.08%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "28.62"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +5.29%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.106"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -1.77%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +0.67%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +1.58%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +2.41%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "53.36"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +2.75%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0₃0737"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +3.52%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0398"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +2.09%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "430.49"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +1.65%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.020.78"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +4.55%  "
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +1.70%  "
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +2.74%  "
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -6.37%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.269"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +1.73%  "
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +4.40%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "26.28"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +1.56%  "
$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = "ThetaToken"
$ws.Range("C49").NumberFormat = "@"
$ws.Range("C49").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.33"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +1.84%  "
$ws.Range("B50").NumberFormat = "@"
$ws.Range("B50").Value = "Stellar"
$ws.Range("C50").NumberFormat = "@"
$ws.Range("C50").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.115"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +1.47%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "35.55"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +14.35%  "
